# MBR Template Map.xlsx - incidental content fixes
#
# Commit: "Update .gitignore to exclude temp Excel files and remove
# temporary web config file" - the workbook-level edits that came along
# with that commit are:
#   1) The "LLIF Dimension Measure" filename in the Dim-Template column
#      (C4) had its stray "_void" suffix removed (the file it points to
#      is no longer voided).
#   2) The sheet's active selection moved from C9 to C11 (just where the
#      author's cursor happened to be when they last saved).
#
# (The raw-XML diff also shows a handful of purely cosmetic / local-machine
# artifacts - the author's absolute folder path, the coauthoring revision
# GUID, the desktop window position, and sub-pixel font-metric noise in
# row heights/column widths/dyDescent - none of which are meaningful
# document content and none of which are exposed anywhere on the Excel
# object model, so there is nothing a script can "do" to reproduce them;
# they are artifacts of the real Excel client on the author's PC.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Fix the filename text in C4: "..._final_void.xlsx" -> "..._final.xlsx"
$ws.Range("C4").Value = "F-825-247M CMD-C.L00.00.C LLIF Dimension Measure Rev 21_final.xlsx"

# 2) Move the selection/active cell to C11 (matches the saved sheetView)
[void]$ws.Range("C11").Select()
